# Apply the "i like this version" predictions refresh:
#  - the data table (A2:H10) is replaced with a new 9-day window
#    (the previous row 2 is dropped, rows shift up by two, and two brand
#    new trailing days are appended with updated count/prcp/snow/wspd)
#  - the now-unused trailing rows 11:12 are removed so the sheet shrinks
#    back down to A1:H10
#  - the active selection moves to J9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2..10, columns A..H (idx, day, date, count, tavg, prcp, snow, wspd)
$data = @(
    @(4, 356, 45728, 58, 43.9, 0, 0, 8.1),
    @(5, 357, 45729, 89, 49.6, 0, 0, 7.5),
    @(6, 358, 45730, 104, 58.5, 0.15, 0, 12.6),
    @(7, 359, 45731, 151, 54.1, 0.25, 0, 18),
    @(1, 360, 45732, 218, 40.5, 0.05, 1, 12.6),
    @(2, 361, 45733, 54, 39.4, 0, 0, 10.5),
    @(3, 362, 45734, 42, 42.3, 0, 0, 11.1),
    @(4, 363, 45735, 58, 43.7, 0, 0, 11),
    @(5, 364, 45736, 197, 44.8, 0.8, 0.1, 11.2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

# Drop the old trailing rows 11 and 12 (data now ends at row 10)
$ws.Rows("11:12").Delete()

# Move the selection to J9, matching the saved view
$ws.Range("J9").Select() | Out-Null
